$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.992.46"
$ws.Range("E2").Value = "'  +0.01%  "
$ws.Range("D3").Value = "'1.909.54"
$ws.Range("E3").Value = "'  +0.36%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'0.7931"
$ws.Range("E5").Value = "'  +6.40%  "
$ws.Range("D6").Value = "'242.01"
$ws.Range("E6").Value = "'  +0.31%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'0.3165"
$ws.Range("E8").Value = "'  +3.12%  "
$ws.Range("D9").Value = "'26.33"
$ws.Range("E9").Value = "'  +3.05%  "
$ws.Range("D10").Value = "'0.06893"
$ws.Range("E10").Value = "'  -0.01%  "
$ws.Range("D11").Value = "'0.07999"
$ws.Range("D12").Value = "'1.906.09"
$ws.Range("E12").Value = "'  +0.22%  "
$ws.Range("D13").Value = "'0.7439"
$ws.Range("E13").Value = "'  -1.42%  "
$ws.Range("D14").Value = "'5.189"
$ws.Range("E14").Value = "'  -1.46%  "
$ws.Range("D15").Value = "'93.07"
$ws.Range("E15").Value = "'  +1.72%  "
$ws.Range("D16").Value = "'30.002.05"
$ws.Range("E16").Value = "'  +0.01%  "
$ws.Range("D17").Value = "'13.94"
$ws.Range("E17").Value = "'  -0.65%  "
$ws.Range("D18").Value = "'5.872"
$ws.Range("E18").Value = "'  -4.67%  "
$ws.Range("D19").Value = "'245.94"
$ws.Range("E19").Value = "'  +3.46%  "
$ws.Range("D20").Value = "'0.000007742"
$ws.Range("E20").Value = "'  -0.11%  "
$ws.Range("D21").Value = "'0.9998"
$ws.Range("E21").Value = "'  +0.00%  "
$ws.Range("D22").Value = "'2.147.12"
$ws.Range("E22").Value = "'  -0.11%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "'  +0.01%  "
$ws.Range("D24").Value = "'6.839"
$ws.Range("E24").Value = "'  -3.68%  "
$ws.Range("D25").Value = "'168.01"
$ws.Range("E25").Value = "'  +1.02%  "
$ws.Range("D26").Value = "'9.231"
$ws.Range("E26").Value = "'  -0.86%  "
$ws.Range("D27").Value = "'0.1400"
$ws.Range("E27").Value = "'  +10.73%  "
$ws.Range("D28").Value = "'18.89"
$ws.Range("E28").Value = "'  +0.49%  "
$ws.Range("D29").Value = "'2.033"
$ws.Range("E29").Value = "'  -1.14%  "
$ws.Range("D30").Value = "'1.365"
$ws.Range("E30").Value = "'  +1.61%  "
$ws.Range("D31").Value = "'1.519"
$ws.Range("E31").Value = "'  -0.34%  "
$ws.Range("D32").Value = "'4.316"
$ws.Range("E32").Value = "'  +0.27%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.05555"
$ws.Range("E33").Value = "'  +2.45%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.086"
$ws.Range("E34").Value = "'  +0.89%  "
$ws.Range("D35").Value = "'1.255"
$ws.Range("E35").Value = "'  -2.32%  "
$ws.Range("D36").Value = "'0.7344"
$ws.Range("E36").Value = "'  -0.50%  "
$ws.Range("E37").Value = "'  -0.19%  "
$ws.Range("D38").Value = "'0.01925"
$ws.Range("E38").Value = "'  -0.97%  "
$ws.Range("D39").Value = "'2.788"
$ws.Range("E39").Value = "'  +0.91%  "
$ws.Range("D40").Value = "'6.144"
$ws.Range("E40").Value = "'  -1.48%  "
$ws.Range("D41").Value = "'0.4418"
$ws.Range("E41").Value = "'  -0.83%  "
$ws.Range("D42").Value = "'72.24"
$ws.Range("E42").Value = "'  -0.54%  "
$ws.Range("D43").Value = "'0.9999"
$ws.Range("E43").Value = "'  +0.03%  "
$ws.Range("D44").Value = "'0.8364"
$ws.Range("E44").Value = "'  +0.64%  "
$ws.Range("D45").Value = "'1.877"
$ws.Range("D46").Value = "'100.52"
$ws.Range("E46").Value = "'  -0.97%  "
$ws.Range("D47").Value = "'7.549"
$ws.Range("E47").Value = "'  -1.68%  "
$ws.Range("D48").Value = "'989.35"
$ws.Range("E48").Value = "'  +8.66%  "
$ws.Range("D49").Value = "'2.054.40"
$ws.Range("E49").Value = "'  -0.15%  "
$ws.Range("D50").Value = "'36.27"
$ws.Range("E50").Value = "'  -0.77%  "
$ws.Range("D51").Value = "'1.479"
$ws.Range("E51").Value = "'  +0.23%  "
